$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# Sheet "展览" (sheet1) - column F updates
$ws1.Range("F2").Value = 170
$ws1.Range("F3").Value = 180
$ws1.Range("F4").Value = 177
$ws1.Range("F5").Value = 4996
$ws1.Range("F7").Value = 45
$ws1.Range("F9").Value = 552
$ws1.Range("F10").Value = 512
$ws1.Range("F11").Value = 1036
$ws1.Range("F12").Value = 18
$ws1.Range("F13").Value = 1388
$ws1.Range("F14").Value = 3659
$ws1.Range("F15").Value = 410
$ws1.Range("F16").Value = 135
$ws1.Range("F17").Value = 118
$ws1.Range("F18").Value = 80
$ws1.Range("F19").Value = 2653
$ws1.Range("F20").Value = 130
$ws1.Range("F21").Value = 16
$ws1.Range("F22").Value = 85
$ws1.Range("F25").Value = 55
$ws1.Range("F26").Value = 128
$ws1.Range("F27").Value = 59
$ws1.Range("F28").Value = 267

# Sheet "全部类型" (sheet4) - column F updates
$ws4.Range("F2").Value = 170
$ws4.Range("F3").Value = 180
$ws4.Range("F4").Value = 177
$ws4.Range("F6").Value = 4996
$ws4.Range("F8").Value = 45
$ws4.Range("F10").Value = 552
$ws4.Range("F11").Value = 512
$ws4.Range("F12").Value = 1036
$ws4.Range("F13").Value = 18
$ws4.Range("F14").Value = 1388
$ws4.Range("F15").Value = 3659
$ws4.Range("F16").Value = 410
$ws4.Range("F17").Value = 135
$ws4.Range("F18").Value = 118
$ws4.Range("F19").Value = 80
$ws4.Range("F20").Value = 2653
$ws4.Range("F21").Value = 130
$ws4.Range("F22").Value = 16
$ws4.Range("F23").Value = 85
$ws4.Range("F26").Value = 55
$ws4.Range("F27").Value = 128
$ws4.Range("F28").Value = 59
$ws4.Range("F29").Value = 267
